# LearningRoadmapTracker.xlsx — Day 9, 10, 11, 12 update
#
# The "Learn networking fundamentals: IP, DNS, TCP/UDP" row (row 6) moves
# from "In Progress" to "Done": the Status column changes, the
# "In Progress?" checkbox is cleared, and the "Done?" checkbox is checked.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "Done"
$ws.Range("D6").Value = "☐"
$ws.Range("E6").Value = "☑"

# Reflect the author's last selection in the saved view.
$ws.Range("C6").Select()
